# Update Leave Card 12/22/2023 10:59 AM
# Shifts the recurring monthly "PERIOD" dates in rows 170-212 forward by one
# month, and fills in the EARNED value (1.25) for rows 175-179 which were
# previously blank placeholders for future periods that are now in the past.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column A: PERIOD dates, shifted forward one month -------------------
$ws.Range("A170").Value = 44957
$ws.Range("A171").Value = 44985
$ws.Range("A172").Value = 45016
$ws.Range("A174").Value = 45046
$ws.Range("A175").Value = 45077
$ws.Range("A176").Value = 45107
$ws.Range("A177").Value = 45138
$ws.Range("A178").Value = 45169
$ws.Range("A179").Value = 45199
$ws.Range("A180").Value = 45230
$ws.Range("A181").Value = 45260
$ws.Range("A182").Value = 45291
$ws.Range("A183").Value = 45322
$ws.Range("A184").Value = 45351
$ws.Range("A185").Value = 45382
$ws.Range("A186").Value = 45412
$ws.Range("A187").Value = 45443
$ws.Range("A188").Value = 45473
$ws.Range("A189").Value = 45504
$ws.Range("A190").Value = 45535
$ws.Range("A191").Value = 45565
$ws.Range("A192").Value = 45596
$ws.Range("A193").Value = 45626
$ws.Range("A194").Value = 45657
$ws.Range("A195").Value = 45688
$ws.Range("A196").Value = 45716
$ws.Range("A197").Value = 45747
$ws.Range("A198").Value = 45777
$ws.Range("A199").Value = 45808
$ws.Range("A200").Value = 45838
$ws.Range("A201").Value = 45869
$ws.Range("A202").Value = 45900
$ws.Range("A203").Value = 45930
$ws.Range("A204").Value = 45961
$ws.Range("A205").Value = 45991
$ws.Range("A206").Value = 46022
$ws.Range("A207").Value = 46053
$ws.Range("A208").Value = 46081
$ws.Range("A209").Value = 46112
$ws.Range("A210").Value = 46142
$ws.Range("A211").Value = 46173
$ws.Range("A212").Value = 46203

# --- Column C: EARNED, now populated for rows that have become current ---
$ws.Range("C175").Value = 1.25
$ws.Range("C176").Value = 1.25
$ws.Range("C177").Value = 1.25
$ws.Range("C178").Value = 1.25
$ws.Range("C179").Value = 1.25
